$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '66.568.06'
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = '3.432.26'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '570.08'
$ws.Range('E5').Value = '  +2.16%  '
Set-TextValue 'D6' '183.07'
$ws.Range('E6').Value = '  +4.06%  '
Set-TextValue 'D7' '0.631'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('D8').Value = '3.429.84'
$ws.Range('E8').Value = '  +2.44%  '
$ws.Range('E9').Value = '  -0.03%  '
Set-TextValue 'D10' '0.172'
$ws.Range('E10').Value = '  +3.25%  '
Set-TextValue 'D11' '0.641'
$ws.Range('E11').Value = '  +1.08%  '
Set-TextValue 'D12' '56.06'
$ws.Range('E12').Value = '  +4.00%  '
$ws.Range('E13').Value = '  +0.71%  '
Set-TextValue 'D14' '9.36'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').Value = '3.982.39'
$ws.Range('E15').Value = '  +2.53%  '
Set-TextValue 'D16' '18.53'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '3.437.33'
$ws.Range('E17').Value = '  +2.68%  '
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '66.656.32'
$ws.Range('E19').Value = '  +2.69%  '
Set-TextValue 'D20' '12.03'
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('E21').Value = '  +2.18%  '
Set-TextValue 'D22' '483.60'
$ws.Range('E22').Value = '  +6.95%  '
Set-TextValue 'D23' '16.28'
$ws.Range('E23').Value = '  +13.93%  '
Set-TextValue 'D24' '5.03'
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('E25').Value = '  +2.08%  '
Set-TextValue 'D26' '88.88'
$ws.Range('E26').Value = '  +2.51%  '
Set-TextValue 'D27' '2.95'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('E28').Value = '  +1.32%  '
Set-TextValue 'D29' '9.04'
$ws.Range('E29').Value = '  +3.41%  '
Set-TextValue 'D30' '31.28'
$ws.Range('E30').Value = '  +1.12%  '
Set-TextValue 'D31' '7.14'
$ws.Range('E31').Value = '  +7.80%  '
Set-TextValue 'D32' '593.61'
$ws.Range('E32').Value = '  +3.72%  '
Set-TextValue 'D33' '11.67'
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('E35').Value = '  +3.24%  '
Set-TextValue 'D36' '1.00'
$ws.Range('E36').Value = '  -0.02%  '
Set-TextValue 'D37' '0.147'
$ws.Range('E37').Value = '  +4.48%  '
$ws.Range('E38').Value = '  -1.62%  '
Set-TextValue 'D39' '0.385'
$ws.Range('E39').Value = '  +4.02%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D40' '36.18'
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0768'
$ws.Range('E41').Value = '  +3.71%  '
$ws.Range('D42').Value = '3.155.51'
$ws.Range('E42').Value = '  +2.88%  '
Set-TextValue 'D43' '2.95'
$ws.Range('E43').Value = '  +5.29%  '
Set-TextValue 'D44' '0.0427'
$ws.Range('E44').Value = '  +2.30%  '
$ws.Range('E45').Value = '  +3.72%  '
Set-TextValue 'D46' '2.76'
$ws.Range('E46').Value = '  +19.70%  '
Set-TextValue 'D47' '3.21'
$ws.Range('E47').Value = '  +0.95%  '
Set-TextValue 'D48' '0.134'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  +6.55%  '
Set-TextValue 'D50' '1.00'
$ws.Range('E50').Value = '  +0.25%  '
Set-TextValue 'D51' '140.71'

Write-Output "applied edits"
